$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper donor cells that already carry the "filled text" style (s=2, border +
# quotePrefix) used throughout the BOM table, so we can restore that
# formatting after writing new values (writing .Value resets a cell's style).
$styleDonor = $ws.Range("B6")   # s=2 style donor

function Restore-Style($range) {
    $styleDonor.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $excel.CutCopyMode = 0
}

# --- Row 3: 100nF capacitor group loses C6, part number changes, quantity drops to 4 ---
$ws.Range("C3").Value = "C1, C2, C3, C5"
Restore-Style $ws.Range("C3")

$ws.Range("E3").Value = "C307331"
Restore-Style $ws.Range("E3")

$ws.Range("F3").Value = 4

# --- Row 7: LED part gets a more specific comment/description and new LibRef/Footprint ---
$ws.Range("A7").Value = "Blue LED 0603"
Restore-Style $ws.Range("A7")

$ws.Range("B7").Value = "Blue 465~475nm 0603 Light Emitting Diodes (LED) ROHS"
Restore-Style $ws.Range("B7")

$ws.Range("D7").Value = "C72041"
Restore-Style $ws.Range("D7")

$ws.Range("E7").Value = "C72041"
Restore-Style $ws.Range("E7")

# --- Rows 13/14: swap the LDO (U2) row and the ATMEGA328P (U3) row, and fix the
#     LDO designator from the placeholder "U?" to the finalized "U2" ---
$ws.Range("A13").Value = "C183880"
Restore-Style $ws.Range("A13")

$ws.Range("B13").Value = "6V - - - Fixed 3.3V SOT-23-3L Linear Voltage Regulators (LDO) ROHS"
Restore-Style $ws.Range("B13")

$ws.Range("C13").Value = "U2"
Restore-Style $ws.Range("C13")

$ws.Range("D13").Value = "FP-SOT23-IPC_C"
Restore-Style $ws.Range("D13")

$ws.Range("E13").Value = "LN1121P332MR-G"
Restore-Style $ws.Range("E13")

$ws.Range("A14").Value = "ATMEGA328P-AU"
Restore-Style $ws.Range("A14")

$ws.Range("B14").Value = "8-bit AVR Microcontroller, 32KB Flash, 1KB EEPROM, 2KB SRAM, 32-pin TQFP, Industrial Grade (-40°C to 85°C)"
Restore-Style $ws.Range("B14")

$ws.Range("C14").Value = "U3"
Restore-Style $ws.Range("C14")

$ws.Range("D14").Value = "32A_L"
Restore-Style $ws.Range("D14")

$ws.Range("E14").Value = "C14877"
Restore-Style $ws.Range("E14")
